$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Build the new abstract text (uses a right single quotation mark, U+2019,
# matching the author's typographic apostrophe in "owners' consumption").
$newAbstract = "Firms in low and middle-income economies often grow slowly, with limited entrepreneurship and innovation. This paper examines whether firm risk aversion prevents risk taking necessary to grow. Because many firms in developing countries are owner-operated, uncertain investments may directly threaten owners" + [char]0x2019 + " consumption. I develop a model of small firm learning which shows that risk aversion can impede product adoption by preventing experimentation with new goods. I test the model using two field experiments in Kenya. Offering shops stock of a new product with insurance that reduces potential losses but lowers expected profits raises stocking by 50\% rejecting risk neutrality. Inducing firms to try stocking the product with a temporary return policy leads to a 70\% increase in stocking after the intervention ends. Firms experiment more when the continuation value of learning increases. The return policy induces risk averse firms uncertain about demand to stock, who overcome uncertainty by learning. These results show that risk aversion can undermine firm learning, preventing enterprises from stocking profitable goods. "

# Replace the Job Market Paper abstract (row 2, column C) with the updated text.
[void]$ws.Range("C2").Select()
$ws.Range("C2").Value = $newAbstract

# The shorter replacement text wraps to fewer lines, so the row height shrinks.
$ws.Rows.Item(2).RowHeight = 391.5

# Leave the selection where the author's last edit left it.
[void]$ws.Range("C3").Select()
